# Apply cell updates per the diff (cryptos list refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '42.384.25'
$ws.Cells.Item(2, 5).Value = '  +0.99%  '
$ws.Cells.Item(3, 4).Value = '2.293.54'
$ws.Cells.Item(3, 5).Value = '  +0.33%  '
$ws.Cells.Item(4, 5).Value = '  -0.06%  '
$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = '@'
$c.Value = '316.34'
$c.Style = 'Normal'
$ws.Cells.Item(5, 5).Value = '  +1.58%  '
$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = '@'
$c.Value = '102.20'
$c.Style = 'Normal'
$ws.Cells.Item(6, 5).Value = '  -3.42%  '
$c = $ws.Cells.Item(7, 4)
$c.NumberFormat = '@'
$c.Value = '0.629'
$c.Style = 'Normal'
$ws.Cells.Item(7, 5).Value = '  +0.43%  '
$ws.Cells.Item(8, 5).Value = '  -0.02%  '
$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = '@'
$c.Value = '0.605'
$c.Style = 'Normal'
$ws.Cells.Item(9, 5).Value = '  -0.40%  '
$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = '@'
$c.Value = '39.55'
$c.Style = 'Normal'
$ws.Cells.Item(10, 5).Value = '  -2.26%  '
$ws.Cells.Item(11, 5).Value = '  -0.46%  '
$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = '@'
$c.Value = '8.41'
$c.Style = 'Normal'
$ws.Cells.Item(12, 5).Value = '  +1.64%  '
$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = '@'
$c.Value = '0.107'
$c.Style = 'Normal'
$ws.Cells.Item(13, 5).Value = '  +0.64%  '
$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = '@'
$c.Value = '0.959'
$c.Style = 'Normal'
$ws.Cells.Item(14, 5).Value = '  -0.71%  '
$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = '@'
$c.Value = '15.18'
$c.Style = 'Normal'
$ws.Cells.Item(15, 5).Value = '  -1.78%  '
$ws.Cells.Item(16, 4).Value = '2.637.72'
$ws.Cells.Item(16, 5).Value = '  +0.24%  '
$ws.Cells.Item(17, 4).Value = '2.290.34'
$ws.Cells.Item(17, 5).Value = '  +1.07%  '
$ws.Cells.Item(18, 4).Value = '42.208.27'
$ws.Cells.Item(18, 5).Value = '  +0.71%  '
$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = '@'
$c.Value = '7.40'
$c.Style = 'Normal'
$ws.Cells.Item(19, 5).Value = '  -1.28%  '
$ws.Cells.Item(20, 5).Value = '  +0.64%  '
$ws.Cells.Item(21, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(21, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = '@'
$c.Value = '12.06'
$c.Style = 'Normal'
$ws.Cells.Item(21, 5).Value = '  +29.70%  '
$ws.Cells.Item(22, 2).Value = 'Litecoin'
$ws.Cells.Item(22, 3).Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = '@'
$c.Value = '73.29'
$c.Style = 'Normal'
$ws.Cells.Item(22, 5).Value = '  -0.06%  '
$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = '@'
$c.Value = '3.54'
$c.Style = 'Normal'
$ws.Cells.Item(23, 5).Value = '  +2.37%  '
$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = '@'
$c.Value = '276.20'
$c.Style = 'Normal'
$ws.Cells.Item(24, 5).Value = '  +7.74%  '
$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = '@'
$c.Value = '2.26'
$c.Style = 'Normal'
$ws.Cells.Item(25, 5).Value = '  -2.11%  '
$ws.Cells.Item(26, 5).Value = '  -0.64%  '
$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = '@'
$c.Value = '10.82'
$c.Style = 'Normal'
$ws.Cells.Item(27, 5).Value = '  -1.18%  '
$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = '@'
$c.Value = '2.37'
$c.Style = 'Normal'
$ws.Cells.Item(28, 5).Value = '  +3.47%  '
$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = '@'
$c.Value = '22.77'
$c.Style = 'Normal'
$ws.Cells.Item(29, 5).Value = '  +0.39%  '
$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = '@'
$c.Value = '37.32'
$c.Style = 'Normal'
$ws.Cells.Item(30, 5).Value = '  +4.61%  '
$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = '@'
$c.Value = '165.94'
$c.Style = 'Normal'
$ws.Cells.Item(31, 5).Value = '  -0.62%  '
$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = '@'
$c.Value = '0.0873'
$c.Style = 'Normal'
$ws.Cells.Item(32, 5).Value = '  -2.06%  '
$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = '@'
$c.Value = '5.97'
$c.Style = 'Normal'
$ws.Cells.Item(33, 5).Value = '  +3.77%  '
$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = '@'
$c.Value = '0.134'
$c.Style = 'Normal'
$ws.Cells.Item(34, 5).Value = '  +3.28%  '
$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = '@'
$c.Value = '2.66'
$c.Style = 'Normal'
$ws.Cells.Item(35, 5).Value = '  -8.66%  '
$c = $ws.Cells.Item(36, 4)
$c.NumberFormat = '@'
$c.Value = '0.117'
$c.Style = 'Normal'
$ws.Cells.Item(36, 5).Value = '  -1.26%  '
$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = '@'
$c.Value = '4.57'
$c.Style = 'Normal'
$ws.Cells.Item(37, 5).Value = '  -0.35%  '
$ws.Cells.Item(38, 5).Value = '  +2.64%  '
$ws.Cells.Item(39, 5).Value = '  +2.04%  '
$ws.Cells.Item(40, 5).Value = '  -1.70%  '
$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = '@'
$c.Value = '1.50'
$c.Style = 'Normal'
$ws.Cells.Item(41, 5).Value = '  +2.06%  '
$ws.Cells.Item(42, 2).Value = 'MultiversX'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = '@'
$c.Value = '69.84'
$c.Style = 'Normal'
$ws.Cells.Item(42, 5).Value = '  -2.33%  '
$ws.Cells.Item(43, 2).Value = 'BitcoinSV'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = '@'
$c.Value = '96.11'
$c.Style = 'Normal'
$ws.Cells.Item(43, 5).Value = '  -1.64%  '
$ws.Cells.Item(44, 2).Value = 'Algorand'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = '@'
$c.Value = '0.225'
$c.Style = 'Normal'
$ws.Cells.Item(44, 5).Value = '  -0.85%  '
$ws.Cells.Item(45, 2).Value = 'FirstDigitalUSD'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = '@'
$c.Value = '1.00'
$c.Style = 'Normal'
$ws.Cells.Item(45, 5).Value = '  -0.25%  '
$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = '@'
$c.Value = '12.01'
$c.Style = 'Normal'
$ws.Cells.Item(46, 5).Value = '  -1.84%  '
$ws.Cells.Item(47, 2).Value = 'Aave'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = '@'
$c.Value = '113.06'
$c.Style = 'Normal'
$ws.Cells.Item(47, 5).Value = '  +0.75%  '
$ws.Cells.Item(48, 2).Value = 'ordi'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = '@'
$c.Value = '79.62'
$c.Style = 'Normal'
$ws.Cells.Item(48, 5).Value = '  +5.45%  '
$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = '@'
$c.Value = '8.97'
$c.Style = 'Normal'
$ws.Cells.Item(49, 5).Value = '  -0.70%  '
$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = '@'
$c.Value = '5.26'
$c.Style = 'Normal'
$ws.Cells.Item(50, 5).Value = '  -0.53%  '
$ws.Cells.Item(51, 4).Value = '1.594.23'
$ws.Cells.Item(51, 5).Value = '  +2.75%  '
